$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.228339177673717
$ws.Range("C2").Value = 2.291705392391212
$ws.Range("D2").Value = -1.315186645652444
$ws.Range("E2").Value = -1.517449058573539
$ws.Range("F2").Value = 0.5798349304998711
$ws.Range("G2").Value = 0.002474916446533804
$ws.Range("H2").Value = 0.08952975383526968
$ws.Range("I2").Value = -0.1695141488638116
$ws.Range("J2").Value = 0.6667485897783219
$ws.Range("K2").Value = 0.5663006126625323
$ws.Range("B3").Value = 2.699937846012982
$ws.Range("C3").Value = -1.06356360176605
$ws.Range("D3").Value = -1.599355112396674
$ws.Range("E3").Value = 0.6549269798811568
$ws.Range("F3").Value = 0.1083820769602576
$ws.Range("G3").Value = 0.1479463999939251
$ws.Range("H3").Value = -0.1030754650337201
$ws.Range("I3").Value = 0.7414305483891583
$ws.Range("J3").Value = 0.6366139045635111
$ws.Range("K3").Value = 0.2823415042266514
$ws.Range("B4").Value = -1.169062813614865
$ws.Range("C4").Value = -1.725460396362032
$ws.Range("D4").Value = 0.8231750485143884
$ws.Range("E4").Value = 0.1690461695641399
$ws.Range("F4").Value = 0.159079239368558
$ws.Range("G4").Value = -0.04021271370678947
$ws.Range("H4").Value = 0.7987820761829015
$ws.Range("I4").Value = 0.68111908156789
$ws.Range("J4").Value = 0.33341537525307
$ws.Range("K4").Value = 0.4838159431822899
$ws.Range("B5").Value = -2.329901016209771
$ws.Range("C5").Value = 0.758971671159695
$ws.Range("D5").Value = 0.4605934776213567
$ws.Range("E5").Value = 0.1213723632826244
$ws.Range("F5").Value = -0.04879673378025057
$ws.Range("G5").Value = 0.8825231537802698
$ws.Range("H5").Value = 0.7159462264959351
$ws.Range("I5").Value = 0.359979758305367
$ws.Range("J5").Value = 0.5293147839288888
$ws.Range("K5").Value = 0.6788178843816368
$ws.Range("B6").Value = 0.6713308423207838
$ws.Range("C6").Value = 0.4081742755916882
$ws.Range("D6").Value = 0.18001876026763
$ws.Range("E6").Value = -0.04386758844314925
$ws.Range("F6").Value = 0.8714298707313746
$ws.Range("G6").Value = 0.7287129257210216
$ws.Range("H6").Value = 0.3687093173527261
$ws.Range("I6").Value = 0.5318167280777406
$ws.Range("J6").Value = 0.6850938218533179
$ws.Range("K6").Value = 0.1817452544735015
$ws.Range("B7").Value = 0.4306022552246968
$ws.Range("C7").Value = 0.1826222998300787
$ws.Range("D7").Value = -0.06593527409716449
$ws.Range("E7").Value = 0.8654478904270196
$ws.Range("F7").Value = 0.7233696755386105
$ws.Range("G7").Value = 0.3583407839823598
$ws.Range("H7").Value = 0.5233394059541168
$ws.Range("I7").Value = 0.6773215135465072
$ws.Range("J7").Value = 0.1731212809405329
$ws.Range("K7").Value = 0.4559179537672726
$ws.Range("B8").Value = 0.2040468712872988
$ws.Range("C8").Value = 0.05424257470204791
$ws.Range("D8").Value = 0.7762535460132365
$ws.Range("E8").Value = 0.6936426716900327
$ws.Range("F8").Value = 0.3669669040122699
$ws.Range("G8").Value = 0.4973992325286816
$ws.Range("H8").Value = 0.6545429505182796
$ws.Range("I8").Value = 0.1595032457919151
$ws.Range("J8").Value = 0.4373971565820968
$ws.Range("K8").Value = 0.2648072220727796
$ws.Range("B9").Value = 0.01192194238165845
$ws.Range("C9").Value = 0.7470490721575034
$ws.Range("D9").Value = 0.6906220718851115
$ws.Range("E9").Value = 0.3481056204007895
$ws.Range("F9").Value = 0.4770681014634773
$ws.Range("G9").Value = 0.6400127950840317
$ws.Range("H9").Value = 0.1428622002996471
$ws.Range("I9").Value = 0.4198133111533466
$ws.Range("J9").Value = 0.2483333385427368
$ws.Range("K9").Value = 0.5550469433309027
$ws.Range("B10").Value = 1.087246812775413
$ws.Range("C10").Value = 0.7659931336690486
$ws.Range("D10").Value = 0.1570123340877904
$ws.Range("E10").Value = 0.5017620140502106
$ws.Range("F10").Value = 0.6493921986403277
$ws.Range("G10").Value = 0.08805140168825606
$ws.Range("H10").Value = 0.4009192202646768
$ws.Range("I10").Value = 0.2334693377276911
$ws.Range("J10").Value = 0.526335232701532
$ws.Range("K10").Value = -0.08123716056912761
$ws.Range("B11").Value = 1.216688005659299
$ws.Range("C11").Value = 0.2041230363001488
$ws.Range("D11").Value = 0.2614787832594023
$ws.Range("E11").Value = 0.6853372797061905
$ws.Range("F11").Value = 0.08600064424866319
$ws.Range("G11").Value = 0.3243509808286466
$ws.Range("H11").Value = 0.2060699433524077
$ws.Range("I11").Value = 0.5000461825038066
$ws.Range("J11").Value = -0.1246344512824602
$ws.Range("K11").Value = 0.5737435035592049
$ws.Range("B12").Value = 0.5164486232236872
$ws.Range("C12").Value = 0.3947675229949266
$ws.Range("D12").Value = 0.506258857889999
$ws.Range("E12").Value = 0.1141948830192304
$ws.Range("F12").Value = 0.3603762920210401
$ws.Range("G12").Value = 0.1685283731481941
$ws.Range("H12").Value = 0.4945038434164454
$ws.Range("I12").Value = -0.1204591223437116
$ws.Range("J12").Value = 0.562558789819434
$ws.Range("B13").Value = 0.6303580141027678
$ws.Range("C13").Value = 0.5909375987643086
$ws.Range("D13").Value = -0.03140548361448672
$ws.Range("E13").Value = 0.3703300348802827
$ws.Range("F13").Value = 0.183893542847739
$ws.Range("G13").Value = 0.4556555084590223
$ws.Range("H13").Value = -0.1363049313708975
$ws.Range("I13").Value = 0.5542245613519331
$ws.Range("B14").Value = 0.9019617852456914
$ws.Range("C14").Value = 0.08568329079670847
$ws.Range("D14").Value = 0.2078658934307159
$ws.Range("E14").Value = 0.2128524933129406
$ws.Range("F14").Value = 0.4913983809139259
$ws.Range("G14").Value = -0.1676325777545246
$ws.Range("H14").Value = 0.5509995116504074
$ws.Range("B15").Value = 0.3325084682008229
$ws.Range("C15").Value = 0.2255217609686054
$ws.Range("D15").Value = 0.1182302542019461
$ws.Range("E15").Value = 0.5236157691624059
$ws.Range("F15").Value = -0.1552444519268073
$ws.Range("G15").Value = 0.5301357081011285
$ws.Range("B16").Value = 0.4651981203848173
$ws.Range("C16").Value = 0.2032544976711596
$ws.Range("D16").Value = 0.4071311908043919
$ws.Range("E16").Value = -0.1401391802749548
$ws.Range("F16").Value = 0.5533747437271186
$ws.Range("B17").Value = 0.3681145747052469
$ws.Range("C17").Value = 0.4206684630523081
$ws.Range("D17").Value = -0.2031120992649178
$ws.Range("E17").Value = 0.5634129094086165
$ws.Range("B18").Value = 0.6745175049177161
$ws.Range("C18").Value = -0.1037940490828814
$ws.Range("D18").Value = 0.4518363671933744
$ws.Range("B19").Value = -0.0597343578434324
$ws.Range("C19").Value = 0.4732568720679752
$ws.Range("B20").Value = 0.7095000033804217
